$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 12, pushing existing rows 12-16 down to 13-17
$ws.Rows.Item(12).Insert()

# Copy the date style (numFmt) from the row below (now row 13, which used to be row 12)
$ws.Range("D12").Value = 44638

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 100112026
$ws.Range("G12").Value = "Haba"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 950
$ws.Range("M12").Value = 925
$ws.Range("N12").Value = "`$/kilo"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 925
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"

# Make sure the date cell style matches the other date cells in column D (style index 2 / numFmt 165)
$ws.Range("D12").NumberFormat = $ws.Range("D13").NumberFormat
